$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror column O's formatting into the new column P (years header row 4
# through the last data row 14), matching how the 2022 column was added.
$ws.Range("O4:O14").Copy()
$ws.Range("P4:P14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New 2022 data values
$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 1
$ws.Range("P6").Value = "-"
$ws.Range("P7").Value = "-"
$ws.Range("P8").Value = "-"
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = "-"
$ws.Range("P12").Value = 1
$ws.Range("P13").Value = "-"
$ws.Range("P14").Value = "-"

[void]$ws.Range("O21:O22").Select()
